$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns
$ws.Range("BC1").Value = "Odd_CS_3-3_HT"
$ws.Range("BD1").Value = "Odd_CS_4-4_HT"

# Copy header style from an existing header cell (BB1) to the new header cells
$ws.Range("BB1").Copy()
$ws.Range("BC1:BD1").PasteSpecial(-4122)  # xlPasteFormats

# New data values
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51

# Updated values in row 2
$ws.Range("G2").Value = 2.15
$ws.Range("I2").Value = 3.3
$ws.Range("J2").Value = 2.77
$ws.Range("L2").Value = 3.65
$ws.Range("U2").Value = 1.6
$ws.Range("V2").Value = 2.22
$ws.Range("W2").Value = 8.5
$ws.Range("X2").Value = 11
$ws.Range("Z2").Value = 21
$ws.Range("AB2").Value = 23
$ws.Range("AD2").Value = 6.4
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 12.5
$ws.Range("AI2").Value = 21
$ws.Range("AK2").Value = 50
$ws.Range("AM2").Value = 27
$ws.Range("AN2").Value = 4.2
$ws.Range("AO2").Value = 11.5
$ws.Range("AP2").Value = 18
$ws.Range("AR2").Value = 70
$ws.Range("AW2").Value = 5.4
$ws.Range("BA2").Value = 90
$ws.Range("BB2").Value = 200
